# Insert a new weekly price record at row 523 for
# "Femacal de La Calera - Cilantro". This pushes the existing rows
# 523..648 down to 524..649 (dimension grows from A1:R648 to A1:R649),
# matching the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above current row 523 (existing rows shift down by one,
# inheriting the neighbouring row's formatting, including the date style on D).
$ws.Rows.Item(523).Insert()

# Populate the newly inserted row 523 with the new record's values.
$ws.Cells.Item(523, 1).Value  = 3                       # Mercado ID
$ws.Cells.Item(523, 2).Value  = "Femacal de La Calera"  # Mercado
$ws.Cells.Item(523, 3).Value  = "Coquimbo"              # Region
$ws.Cells.Item(523, 4).Value  = 45204                   # Fecha (serial date)
$ws.Cells.Item(523, 5).Value  = 5                       # Codreg
$ws.Cells.Item(523, 6).Value  = 100112040               # Categoria ID
$ws.Cells.Item(523, 7).Value  = "Cilantro"              # Categoria
$ws.Cells.Item(523, 8).Value  = "Sin especificar"       # Variedad
$ws.Cells.Item(523, 9).Value  = "Primera"               # Calidad
$ws.Cells.Item(523, 10).Value = 110                     # Volumen
$ws.Cells.Item(523, 11).Value = 4000                    # Precio minimo
$ws.Cells.Item(523, 12).Value = 4000                    # Precio maximo
$ws.Cells.Item(523, 13).Value = 4000                    # Precio promedio ponderado
$ws.Cells.Item(523, 14).Value = "`$/docena de atados (3 kilos)"  # Unidad de comercializacion
$ws.Cells.Item(523, 15).Value = "Provincia de Quillota"  # Origen
$ws.Cells.Item(523, 16).Value = 1333                     # Precio $/Kg
$ws.Cells.Item(523, 17).Value = 3                        # Kg o Unidades
$ws.Cells.Item(523, 18).Value = "Hortaliza"               # Clasificacion

# Make sure the date cell keeps the same date formatting used throughout column D.
$ws.Cells.Item(523, 4).NumberFormat = $ws.Cells.Item(524, 4).NumberFormat
